# Update gh-pages output (丽水-漫展信息.xlsx) to the newer scrape.
#
# A brand-new exhibition ("丽水·LZ栗子动漫游戏嘉年华", 2024-07-14) was added as the
# new row 2 on both the "展览" and "全部类型" sheets, pushing every later row down
# by one. The running index in column A is renumbered to stay sequential, and two
# "want to go" counters (F3 and F9) ticked up since the previous snapshot.

function Update-LsSheet {
    param($ws)

    # Insert a fresh row right under the header; Excel shifts rows 2-8 down to 3-9
    # and keeps their values/formats intact.
    $ws.Rows.Item(2).Insert()

    # The new row needs the same bold/centered/bordered look as the rest of
    # column A, so copy that formatting from the row right below (now row 3).
    $ws.Range("A3").Copy()
    $ws.Range("A2").PasteSpecial(-4122)  # xlPasteFormats
    $ws.Range("A2").Value = 1

    # Column B holds plain text dates ("2024-07-14"); force text formatting
    # first so Excel doesn't auto-convert the literal into a date serial,
    # then drop the temporary number format again so the cell matches the
    # plain/default styling used by every other row.
    $ws.Range("B2").NumberFormat = "@"
    $ws.Range("B2").Value = "2024-07-14"
    $ws.Range("B2").ClearFormats()

    $ws.Range("C2").Value = "丽水·LZ栗子动漫游戏嘉年华"
    $ws.Range("D2").Value = "城北街798号 莱茵体育生活馆"
    $ws.Range("E2").Value = "2024.07.14 09:30-07.14 17:00"
    $ws.Range("F2").Value = 0
    $ws.Range("G2").Value = 50
    $ws.Range("H2").Value = "https://show.bilibili.com/platform/detail.html?id=87480"
    $ws.Range("I2").Value = "//i2.hdslb.com/bfs/openplatform/202406/5F9bgOvv1718607603808.jpeg"

    # Renumber the running index in column A for every shifted row (3..9 -> 2..8).
    for ($r = 3; $r -le 9; $r++) {
        $ws.Cells.Item($r, 1).Value = $r - 1
    }

    # "Want to go" counters that ticked up between scrapes.
    $ws.Range("F3").Value = 74   # 丽水·CCAC动漫游戏嘉年华: 73 -> 74
    $ws.Range("F9").Value = 227  # 丽水·AEO纯白礼赞动漫嘉年华: 224 -> 227
}

$wb = $excel.ActiveWorkbook

Update-LsSheet($wb.Worksheets.Item("展览"))
Update-LsSheet($wb.Worksheets.Item("全部类型"))
